$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper staging cell used to write literal "True"/"False" text without Excel
# auto-converting the input to a Boolean value (PasteSpecial xlPasteValues = -4163
# copies the computed *text* result of a ="..."  formula, landing as a real string cell).
$stage = $ws.Range("ZZ1")

function Set-TextValue($cellRange, [string]$text) {
    $stage.Formula = ('="' + $text + '"')
    $stage.Copy()
    $cellRange.PasteSpecial(-4163)
}

# Update dimension-driving data rows 2-22 (A:E)
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 'Aircraft ActiveTrack available at max speed . Obstacle Avoidance is not available .'
$ws.Range("C2").Value = 'Aircraft ActiveTrack available at max speed'
$ws.Range("D2").Value = '0-5'
$ws.Range("E2").Value = 'Missing'

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 'Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available .'
$ws.Range("C3").Value = 'Aircraft ActiveTrack available at max speed'
$ws.Range("D3").Value = '0-5'
$ws.Range("E3").Value = 'Missing'

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 'Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available .'
$ws.Range("C4").Value = 'When exceeding nnn, Obstacle Avoidance is not available'
$ws.Range("D4").Value = '7-14'
$ws.Range("E4").Value = 'Missing'

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 'Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available .'
$ws.Range("C5").Value = 'When exceeding nnn,'
$ws.Range("D5").Value = '7-9'
Set-TextValue $ws.Range("E5") 'False'

$ws.Range("A6").Value = 48
$ws.Range("B6").Value = 'Check whether propellers are installed correctly . If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist . Contact DJI Support for assistance .'
$ws.Range("C6").Value = 'If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist'
$ws.Range("D6").Value = '7-23'
$ws.Range("E6").Value = 'Missing'

$ws.Range("A7").Value = 48
$ws.Range("B7").Value = 'Check whether propellers are installed correctly . If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist . Contact DJI Support for assistance .'
$ws.Range("C7").Value = 'If the propellers are installed correctly and the aircraft still cannot takeoff,'
$ws.Range("D7").Value = '7-18'
Set-TextValue $ws.Range("E7") 'False'

$ws.Range("A8").Value = 50
$ws.Range("B8").Value = 'Compass abnormal . Solution: 1. Ensure there are no magnets or metal objects near the aircraft . The ground or walls may contain metal . Move away from sources of interference before attempting flight . 2. Calibrate Compass Before Takeoff .'
$ws.Range("C8").Value = '2. Calibrate Compass Before Takeoff'
$ws.Range("D8").Value = '35-39'
$ws.Range("E8").Value = 'Missing'

$ws.Range("A9").Value = 50
$ws.Range("B9").Value = 'Compass abnormal . Solution: 1. Ensure there are no magnets or metal objects near the aircraft . The ground or walls may contain metal . Move away from sources of interference before attempting flight . 2. Calibrate Compass Before Takeoff .'
$ws.Range("C9").Value = 'Calibrate Compass Before Takeoff'
$ws.Range("D9").Value = '36-39'
Set-TextValue $ws.Range("E9") 'False'

$ws.Range("A10").Value = 66
$ws.Range("B10").Value = 'Downlink data connection lost for nnn seconds .'
$ws.Range("C10").Value = 'Downlink data connection lost for nnn seconds'
$ws.Range("D10").Value = '0-6'
$ws.Range("E10").Value = 'Missing'

$ws.Range("A11").Value = 66
$ws.Range("B11").Value = 'Downlink data connection lost for nnn seconds .'
$ws.Range("C11").Value = 'Downlink data connection lost for nnn'
$ws.Range("D11").Value = '0-5'
Set-TextValue $ws.Range("E11") 'False'

$ws.Range("A12").Value = 70
$ws.Range("B12").Value = 'Downward ambient light too low . Obstacle avoidance unavailable . Fly with caution . Backward ambient light too low . Backward obstacle avoidance unavailable . Only infrared sensors available . Fly with caution .'
$ws.Range("C12").Value = 'Backward ambient light too low'
$ws.Range("D12").Value = '14-18'
$ws.Range("E12").Value = 'Missing'

$ws.Range("A13").Value = 77
$ws.Range("B13").Value = 'Exiting GPS mode : Unknown Error .'
$ws.Range("C13").Value = 'Unknown Error'
$ws.Range("D13").Value = '4-5'
$ws.Range("E13").Value = 'Missing'

$ws.Range("A14").Value = 81
$ws.Range("B14").Value = 'Extra payload detected . Return aircraft to an area nearby the home point promptly and fly in a wind-free environment to ensure flight safety .'
$ws.Range("C14").Value = 'Return aircraft to an area nearby the home point promptly and fly in a wind-free environment to ensure flight safety'
$ws.Range("D14").Value = '4-23'
$ws.Range("E14").Value = 'Missing'

$ws.Range("A15").Value = 86
$ws.Range("B15").Value = 'Flight altitude exceeds nnn . May violate local policies and regulations . Ensure you have obtained proper airspace authorization .'
$ws.Range("C15").Value = 'May violate local policies and regulations'
$ws.Range("D15").Value = '5-10'
$ws.Range("E15").Value = 'Missing'

$ws.Range("A16").Value = 91
$ws.Range("B16").Value = 'GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn .'
$ws.Range("C16").Value = 'GEO Zone Info: The target area is in an Altitude Zone'
$ws.Range("D16").Value = '0-10'
$ws.Range("E16").Value = 'Missing'

$ws.Range("A17").Value = 91
$ws.Range("B17").Value = 'GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn .'
$ws.Range("C17").Value = 'GEO Zone Info:'
$ws.Range("D17").Value = '0-2'
Set-TextValue $ws.Range("E17") 'False'

$ws.Range("A18").Value = 91
$ws.Range("B18").Value = 'GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn .'
$ws.Range("C18").Value = 'The target area is in an Altitude Zone'
$ws.Range("D18").Value = '3-10'
Set-TextValue $ws.Range("E18") 'False'

$ws.Range("A19").Value = 113
$ws.Range("B19").Value = 'Insufficient SD card space . Change card or delete files .'
$ws.Range("C19").Value = 'Insufficient SD card space'
$ws.Range("D19").Value = '0-3'
$ws.Range("E19").Value = 'Missing'

$ws.Range("A20").Value = 115
$ws.Range("B20").Value = 'Landin .'
$ws.Range("C20").Value = 'Landin'
$ws.Range("D20").Value = '0-0'
$ws.Range("E20").Value = 'Missing'

$ws.Range("A21").Value = 147
$ws.Range("B21").Value = 'SD card write speed is too slow . Not suitable for shooting a 4K video .'
$ws.Range("C21").Value = 'Not suitable for shooting a 4K video'
$ws.Range("D21").Value = '8-14'
$ws.Range("E21").Value = 'Missing'

$ws.Range("A22").Value = 147
$ws.Range("B22").Value = 'SD card write speed is too slow . Not suitable for shooting a 4K video .'
$ws.Range("C22").Value = 'Not suitable for shooting a 4K'
$ws.Range("D22").Value = '8-13'
Set-TextValue $ws.Range("E22") 'False'

# Clean up the staging cell so it does not show up as stray data
$stage.Clear()
$excel.CutCopyMode = $false